# Updated cryptos list cell values (Price + Volume(1h)) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.015.21'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.643.95'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").Formula = "'216.56"
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("D9").Formula = "'0.0641"
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("D10").Formula = "'19.67"
$ws.Range("E10").Value = '  +0.27%  '
$ws.Range("D11").Formula = "'0.0795"
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Formula = "'4.30"
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("D13").Value = '1.869.89'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '1.633.29'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").Formula = "'63.11"
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").Value = '25.988.29'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Formula = "'193.36"
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").Value = '  +7.30%  '
$ws.Range("E25").Value = '  +1.93%  '
$ws.Range("D26").Formula = "'144.66"
$ws.Range("E26").Value = '  +1.36%  '
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("D29").Formula = "'15.55"
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("E31").Value = '  +0.66%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").Formula = "'3.26"
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = '1.132.99'
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("E38").Value = '  -0.93%  '
$ws.Range("D39").Formula = "'2.47"
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").Formula = "'5.52"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").Formula = "'99.31"
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '1.778.65'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").Value = '0.0₆0116'
$ws.Range("E45").Value = '  +4.33%  '
$ws.Range("D46").Formula = "'56.70"
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").Formula = "'7.74"
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("D51").Formula = "'0.0961"
$ws.Range("E51").Value = '  +0.39%  '
